$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - updated (was old row 3 data, F value changed)
$ws.Range("A1").Value = 1.0
$ws.Range("B1").Value = "MX-t11"
$ws.Range("C1").Value = "Panadol"
$ws.Range("D1").Value = 1.9
$ws.Range("E1").Value = 14.0
$ws.Range("F1").Value = 2570.0
$ws.Range("G1").Value = "Yellowish"
$ws.Range("H1").Value = "Nein"

# Row 2 - updated (new medicine names, was old row 1 red/shape, F value changed)
$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = "Abxx"
$ws.Range("C2").Value = "demoMednew"
$ws.Range("D2").Value = 1.5
$ws.Range("E2").Value = 10.0
$ws.Range("F2").Value = 310.0
$ws.Range("G2").Value = "red"
$ws.Range("H2").Value = "shape"

# Row 3 - new medicine
$ws.Range("A3").Value = -1.0
$ws.Range("B3").Value = "ab"
$ws.Range("C3").Value = "mai thuy"
$ws.Range("D3").Value = 1.0
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 40.0
$ws.Range("G3").Value = "do"
$ws.Range("H3").Value = "tron"

# Row 4 - new medicine
$ws.Range("A4").Value = -1.0
$ws.Range("B4").Value = "hmm"
$ws.Range("C4").Value = "thuoc te"
$ws.Range("D4").Value = 1.0
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 400.0
$ws.Range("G4").Value = "tim"
$ws.Range("H4").Value = "vuong"
